# Audi.xlsx / CarDetails sheet — trim the test fixture down to a smaller
# column set (keep Registration Number, Vehicle make, Vehicle colour,
# Date of first registration, Year of manufacture, Fuel type; drop
# Cylinder capacity, CO2 Emissions, Euro Status, Export marker,
# Vehicle status, Vehicle type approval, Wheelplan, Revenue weight).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear formatting first so the column-width/bestFit metadata for the
# removed columns is dropped cleanly instead of leaving stale <col> entries
# behind once the columns themselves are deleted.
$ws.Columns("F:G").ClearFormats()
$ws.Columns("I:N").ClearFormats()

# Delete the unwanted columns, right-to-left so earlier deletes don't
# invalidate the letter-ranges of later ones.
$ws.Columns("I:N").Delete()
$ws.Columns("F:G").Delete()

# Match the view state recorded in the saved workbook: scrolled right so
# column E is left-most on screen, with I9 as the active (empty) selection.
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("I9").Select()
